$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

$ws.Range("N1").Value = "description"

$ws.Range("N1").Select()
